# Update column G ("K") values on Sheet1 to reflect the regenerated
# save_data (K replaces old Strike# values; std/mean recalculated upstream).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newValues = @{
    2  = 9
    3  = 2
    4  = 2
    5  = 9
    6  = 6
    7  = 8
    8  = 1
    9  = 6
    10 = 4
    11 = 8
    12 = 5
    13 = 8
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
